# Update the dSF (F column) values for the specified rows to reflect the
# repulled data / recalculated mean values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    7  = 7
    9  = -1
    10 = -8
    11 = -6
    12 = -1
    16 = -2
    26 = -6
    28 = 3
    29 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
